$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Sema3c-Nrp1 results (Dr Hou advice): sending clusters now
# include ECs, FAPs and sCs (previously only FAPs and sCs), each paired
# with target clusters ECs, FAPs, M2 and sCs -> 12 data rows total.
$data = @(
    @("ECs", "Sema3c", "Nrp1", "ECs", 2, 0.6666666666666666, 0.2583983333333333, 0.775195, 0.007195239230717037, 0.007195239230717038, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 30.50315124191278, 274.5283611772149, 0.003350978507054935, 0.003350978507054935),
    @("ECs", "Sema3c", "Nrp1", "FAPs", 2, 0.6666666666666666, 0.2583983333333333, 0.775195, 0.007195239230717037, 0.007195239230717038, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 15.34004065902444, 138.06036593122, 0.001685207739294424, 0.001685207739294424),
    @("ECs", "Sema3c", "Nrp1", "M2", 2, 0.6666666666666666, 0.2583983333333333, 0.775195, 0.007195239230717037, 0.007195239230717038, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 14.68866988195945, 132.198028937635, 0.001613650231784534, 0.001613650231784534),
    @("ECs", "Sema3c", "Nrp1", "sCs", 2, 0.6666666666666666, 0.2583983333333333, 0.775195, 0.007195239230717037, 0.007195239230717038, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 4.964670055261111, 44.68203049735, 0.0005454027525831454, 0.0005454027525831455),
    @("FAPs", "Sema3c", "Nrp1", "ECs", 3, 1, 35.10542566666667, 105.316277, 0.9775292770250872, 0.9775292770250873, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 4144.090616639916, 37296.81554975925, 0.4552565234167454, 0.4552565234167456),
    @("FAPs", "Sema3c", "Nrp1", "FAPs", 3, 1, 35.10542566666667, 105.316277, 0.9775292770250872, 0.9775292770250873, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 2084.063972596677, 18756.57575337009, 0.228948593675237, 0.228948593675237),
    @("FAPs", "Sema3c", "Nrp1", "M2", 3, 1, 35.10542566666667, 105.316277, 0.9775292770250872, 0.9775292770250873, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 1995.570180470718, 17960.13162423646, 0.219226949079566, 0.2192269490795661),
    @("FAPs", "Sema3c", "Nrp1", "sCs", 3, 1, 35.10542566666667, 105.316277, 0.9775292770250872, 0.9775292770250873, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 674.4890856539123, 6070.40177088521, 0.07409721085353879, 0.07409721085353882),
    @("sCs", "Sema3c", "Nrp1", "ECs", 3, 1, 0.5485793333333333, 1.645738, 0.01527548374419571, 0.01527548374419571, 3, 1, 118.0470123333333, 354.141037, 0.4657216250363638, 0.4657216250363638, 64.75815132781177, 582.823361950306, 0.007114123112563386, 0.007114123112563387),
    @("sCs", "Sema3c", "Nrp1", "FAPs", 3, 1, 0.5485793333333333, 1.645738, 0.01527548374419571, 0.01527548374419571, 3, 1, 59.36586533333332, 178.097596, 0.2342114953037475, 0.2342114953037476, 32.56688682731644, 293.101981445848, 0.003577693889216166, 0.003577693889216168),
    @("sCs", "Sema3c", "Nrp1", "M2", 3, 1, 0.5485793333333333, 1.645738, 0.01527548374419571, 0.01527548374419571, 3, 1, 56.84506433333333, 170.535193, 0.2242663767030476, 0.2242663767030477, 31.18402749527044, 280.6562474574339, 0.003425777391697076, 0.003425777391697077),
    @("sCs", "Sema3c", "Nrp1", "sCs", 3, 1, 0.5485793333333333, 1.645738, 0.01527548374419571, 0.01527548374419571, 3, 1, 19.21324333333333, 57.63973, 0.07580050295684103, 0.07580050295684104, 10.53998821897111, 94.85989397073999, 0.001157889350719084, 0.001157889350719084)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
